# Update cryptocurrency price/volume data in the worksheet to reflect
# the latest scrape from coinranking.com (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.793.55"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.638.85"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'215.57"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.258"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D11").Value = "'0.0793"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.865.39"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "1.631.36"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'0.562"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "'63.12"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "25.834.45"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'4.47"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").Value = "'192.31"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "'9.96"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "'6.31"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'1.86"
$ws.Range("E24").Value = "  +5.86%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'141.88"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").Value = "'6.92"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "1.134.00"
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("D38").Value = "'2.53"
$ws.Range("D39").Value = "'0.545"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "'5.57"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "'100.56"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "'0.805"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "1.774.31"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +4.02%  "
$ws.Range("D47").Value = "'55.26"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "'1.41"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.31"
$ws.Range("E51").Value = "  -3.25%  "
